$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10: rename header labels and add two new columns (J10, L10)
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = "Cloth"
$ws.Range("F10").Value = "Hoody"
$ws.Range("H10").Value = "Paper"
$ws.Range("J10").Value = "Sweater"
$ws.Range("L10").Value = "Tshirt"

# ---------------------------------------------------------------------------
# Row 11: add Depth/Normal sub-headers for the two new columns
# ---------------------------------------------------------------------------
$ws.Range("J11").Value = "Depth"
$ws.Range("K11").Value = "Normal"
$ws.Range("L11").Value = "Depth"
$ws.Range("M11").Value = "Normal"

# ---------------------------------------------------------------------------
# Clear the old placeholder cells (A12, A13, A14) before building the new table
# ---------------------------------------------------------------------------
$ws.Range("A12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()

# ---------------------------------------------------------------------------
# Row 12 - Baseline / Fixed
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Baseline"
$ws.Range("A3").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B12").Formula = "=231624+268281+268299"
$ws.Range("B3").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C12").Value = "Fixed"
$ws.Range("D12").Value = 0.037727453000000001
$ws.Range("E12").Value = 1.0407907999999999
$ws.Range("F12").Value = 0.088872690000000004
$ws.Range("G12").Value = 1.298073
$ws.Range("H12").Value = 0.064535910000000002
$ws.Range("I12").Value = 1.3023244
$ws.Range("J12").Value = 0.078107949999999995
$ws.Range("K12").Value = 1.3971530000000001
$ws.Range("L12").Value = 0.062479510000000002
$ws.Range("M12").Value = 1.1054889000000001

# ---------------------------------------------------------------------------
# Row 13 - Baseline / Variable
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "Variable"
$ws.Range("D13").Value = 0.036929759999999999
$ws.Range("E13").Value = 1.0440722
$ws.Range("F13").Value = 0.093535660000000007
$ws.Range("G13").Value = 1.3148359999999999
$ws.Range("H13").Value = 0.057419900000000003
$ws.Range("I13").Value = 1.3108557000000001
$ws.Range("J13").Value = 0.070308579999999996
$ws.Range("K13").Value = 1.4024817000000001
$ws.Range("L13").Value = 0.060015798000000002
$ws.Range("M13").Value = 1.0940189
$ws.Range("N13").Value = "<--------"
$ws.Range("O13").Value = "After 10 epochs in this table :)"
$ws.Range("A4").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 14 - UNet / Fixed
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "UNet"
$ws.Range("A5").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B14").Value = 1989452
$ws.Range("B5").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C14").Value = "Fixed"
$ws.Range("D14").Value = 0.033315957
$ws.Range("E14").Value = 0.88682019999999995
$ws.Range("F14").Value = 0.070560250000000005
$ws.Range("G14").Value = 1.1872355000000001
$ws.Range("H14").Value = 0.047023280000000001
$ws.Range("I14").Value = 1.1152321000000001
$ws.Range("J14").Value = 0.050319959999999997
$ws.Range("K14").Value = 1.190321
$ws.Range("L14").Value = 0.049251344000000002
$ws.Range("M14").Value = 0.97129509999999997

# ---------------------------------------------------------------------------
# Row 15 - UNet / Variable
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = "Variable"
$ws.Range("D15").Value = 0.051807974
$ws.Range("E15").Value = 1.5437989000000001
$ws.Range("F15").Value = 0.084237350000000003
$ws.Range("G15").Value = 1.9795997000000001
$ws.Range("H15").Value = 0.063675910000000002
$ws.Range("I15").Value = 1.8204632000000001
$ws.Range("J15").Value = 0.076809310000000006
$ws.Range("K15").Value = 2.0459866999999998
$ws.Range("L15").Value = 0.06038789
$ws.Range("M15").Value = 1.8239620000000001
$ws.Range("A6").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 16 - InceptionNet / Fixed
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "InceptionNet"
$ws.Range("A7").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("B16").Value = 10756132
$ws.Range("B7").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Fixed"
$ws.Range("D16").Value = 0.063842040000000003
$ws.Range("E16").Value = 1.5262985
$ws.Range("F16").Value = 0.10952318
$ws.Range("G16").Value = 1.9895621999999999
$ws.Range("H16").Value = 0.092206389999999999
$ws.Range("I16").Value = 1.8147032000000001
$ws.Range("J16").Value = 0.094131770000000003
$ws.Range("K16").Value = 2.0846285999999998
$ws.Range("L16").Value = 0.081547819999999993
$ws.Range("M16").Value = 1.8472519000000001

# ---------------------------------------------------------------------------
# Row 17 - InceptionNet / Variable
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = "Variable"
$ws.Range("D17").Value = 0.056296843999999999
$ws.Range("E17").Value = 1.5438864999999999
$ws.Range("F17").Value = 0.107318565
$ws.Range("G17").Value = 1.9878129
$ws.Range("H17").Value = 0.090307854000000007
$ws.Range("I17").Value = 1.8178124
$ws.Range("J17").Value = 0.08854476
$ws.Range("K17").Value = 2.0514804999999998
$ws.Range("L17").Value = 0.079266009999999998
$ws.Range("M17").Value = 1.8279300000000001
$ws.Range("A8").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Merge the A/B column pairs for each new net block, like the first table
# ---------------------------------------------------------------------------
$ws.Range("A12:A13").Merge()
$ws.Range("B12:B13").Merge()
$ws.Range("A14:A15").Merge()
$ws.Range("B14:B15").Merge()
$ws.Range("A16:A17").Merge()
$ws.Range("B16:B17").Merge()

# ---------------------------------------------------------------------------
# Update selection to mirror the final cursor position in the authored file
# ---------------------------------------------------------------------------
$ws.Range("M15").Select()
